# Added MyConstructor for changing typical parameters
#
# The table on Sheet1 lists, per cortical layer, where the apical
# dendrites of pyramidal cells (PC1 / PC2) project to, plus a comment.
# L4CA and L4CB previously had placeholder "--" entries; this edit fills
# them in properly and adds a new "L4C" row (a combined/generic L4C
# layer) right after L4CB and before L5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 8 (which holds "L5"), shifting
# L5 and L6 down to rows 9 and 10.
$ws.Rows.Item(8).Insert()

# --- Row 6: L4CA ------------------------------------------------------
$ws.Range("D6").Value = "No pyramidal cells in monkeys, but here for technical reasons/other species"
$ws.Range("B6").Value = "[L4CA->L1]"
$ws.Range("C6").Value = "[L4CA->L23]"

# --- Row 7: L4CB ------------------------------------------------------
$ws.Range("B7").Value = "[L4CB->L1]"
$ws.Range("C7").Value = "[L4CB->L23]"

# --- Row 8 (new): L4C --------------------------------------------------
$ws.Range("A8").Value = "L4C"
$ws.Range("B8").Value = "[L4C->L1]"
$ws.Range("C8").Value = "[L4C->L23]"
$ws.Range("D8").Value = "No pyramidal cells"

# Selection bookkeeping, mirroring the saved state in the diff.
$ws.Range("D16").Select()
